$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the title font used by B2 (Arial): drop italic + red color, shrink to 18pt ---
$titleFont = $ws.Range("B2").Font
$titleFont.Italic = $false
$titleFont.Size = 18
$titleFont.ColorIndex = 1

# --- Replace the old sample data rows with a description/header row ---
# Remove row 3 ("32","33","34") and row 5 ("52","53","54") entirely.
$ws.Rows("3:3").ClearContents()
$ws.Rows("5:5").ClearContents()

# Turn row 4 into the new header row describing the object/title columns.
$ws.Range("B4").Value = "type"
$ws.Range("C4").Value = "object"
$ws.Range("D4").Value = "subsystem"
$ws.Range("E4").Value = "author"
$ws.Range("F4").Value = "insert"
$ws.Range("G4").Value = "delete"
